$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts existing rows 2..320 down to 3..321)
$ws.Rows.Item(2).Insert()

# Force the new row's cells to Text format so numeric-looking / date-looking
# strings are not auto-converted to numbers/dates (keeps literal content,
# matching how the rest of the sheet stores its values as plain text).
$ws.Range("A2:C2").NumberFormat = "@"

# Populate the newly inserted row with the new record
$ws.Cells.Item(2, 1).Value = "+5521994075389"
$ws.Cells.Item(2, 2).Value = "21"
$ws.Cells.Item(2, 3).Value = "2024-10-16"

# Re-apply the same cell style used by the rest of the data rows (copy
# formats only from the row immediately below, which holds the data that
# used to be row 2 before the insert) so the new row matches the sheet's
# existing look instead of inheriting the header row's style.
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
